$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "SIS Login ID" column (column B) is being removed from the Canvas
# course-data export so the remaining columns match Canvas' column naming.
# Select the whole column first (as a user would before right-click >
# Delete), then delete it. This shifts everything after it one column to
# the left (C->B, D->C, ... BK->BJ) and keeps all the row data intact.
$ws.Columns.Item(2).Select() | Out-Null
$ws.Columns.Item(2).Delete() | Out-Null

$wb.Save()
